# The crawl was re-run, so the "timestamp" column (O) on every data row
# needs to reflect the new crawl time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestampCol = 15   # column O
$newTimestamp = "2023-03-10 12:57:21"

# Find the last used row based on column O (xlUp from the bottom of the sheet).
$lastRow = $ws.Cells.Item($ws.Rows.Count, $timestampCol).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $timestampCol).Value = $newTimestamp
}
